$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 2 de Abril de 2020 a las 11:20"

# --- Refresh the country ranking table for the 11:20 data pull. ---
# A few countries swapped rank order because their case counts crossed
# over; for each affected row we set both the country label (column A)
# and the updated stats (columns B:H).

$ws.Range("A14").Value = "Belgica"
$ws.Range("B14").Value = 15348
$ws.Range("C14").Value = 1384
$ws.Range("D14").Value = 2495
$ws.Range("E14").Value = 11842
$ws.Range("F14").Value = 1144
$ws.Range("G14").Value = 183
$ws.Range("H14").Value = 1011

$ws.Range("A16").Value = "Austria"
$ws.Range("B16").Value = 10842
$ws.Range("C16").Value = 131
$ws.Range("D16").Value = 1749
$ws.Range("E16").Value = 8935
$ws.Range("F16").Value = 227
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = 158

$ws.Range("A29").Value = "Malasia"
$ws.Range("B29").Value = 3116
$ws.Range("C29").Value = 208
$ws.Range("D29").Value = 767
$ws.Range("E29").Value = 2299
$ws.Range("F29").Value = 105
$ws.Range("G29").Value = 5
$ws.Range("H29").Value = 50

$ws.Range("A30").Value = "Chile"
$ws.Range("B30").Value = 3031
$ws.Range("C30").Value = 0
$ws.Range("D30").Value = 234
$ws.Range("E30").Value = 2781
$ws.Range("F30").Value = 31
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 16

$ws.Range("A32").Value = "Polonia"
$ws.Range("B32").Value = 2633
$ws.Range("C32").Value = 79
$ws.Range("D32").Value = 56
$ws.Range("E32").Value = 2532
$ws.Range("F32").Value = 50
$ws.Range("G32").Value = 2
$ws.Range("H32").Value = 45

$ws.Range("A33").Value = "Filipinas"
$ws.Range("B33").Value = 2633
$ws.Range("C33").Value = 322
$ws.Range("D33").Value = 51
$ws.Range("E33").Value = 2475
$ws.Range("F33").Value = 1
$ws.Range("G33").Value = 11
$ws.Range("H33").Value = 107

$ws.Range("A34").Value = "Rumania"
$ws.Range("B34").Value = 2460
$ws.Range("C34").Value = 0
$ws.Range("D34").Value = 252
$ws.Range("E34").Value = 2114
$ws.Range("F34").Value = 57
$ws.Range("G34").Value = 2
$ws.Range("H34").Value = 94

$ws.Range("A35").Value = "Japon"
$ws.Range("B35").Value = 2384
$ws.Range("C35").Value = 0
$ws.Range("D35").Value = 472
$ws.Range("E35").Value = 1855
$ws.Range("F35").Value = 69
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 57

$ws.Range("A36").Value = "Luxemburgo"
$ws.Range("B36").Value = 2319
$ws.Range("C36").Value = 0
$ws.Range("D36").Value = 80
$ws.Range("E36").Value = 2210
$ws.Range("F36").Value = 31
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = 29

$ws.Range("A61").Value = "Hong Kong"
$ws.Range("B61").Value = 802
$ws.Range("C61").Value = 36
$ws.Range("D61").Value = 154
$ws.Range("E61").Value = 644
$ws.Range("F61").Value = 5
$ws.Range("G61").Value = 0
$ws.Range("H61").Value = 4

$ws.Range("A62").Value = "Nueva Zelanda"
$ws.Range("B62").Value = 797
$ws.Range("C62").Value = 89
$ws.Range("D62").Value = 92
$ws.Range("E62").Value = 704
$ws.Range("F62").Value = 2
$ws.Range("G62").Value = 0
$ws.Range("H62").Value = 1

$ws.Range("A63").Value = "Egipto"
$ws.Range("B63").Value = 779
$ws.Range("C63").Value = 0
$ws.Range("D63").Value = 179
$ws.Range("E63").Value = 548
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 52

$ws.Range("A95").Value = "Vietnam"
$ws.Range("B95").Value = 222
$ws.Range("C95").Value = 4
$ws.Range("D95").Value = 75
$ws.Range("E95").Value = 147
$ws.Range("F95").Value = 3
$ws.Range("G95").Value = 0
$ws.Range("H95").Value = 0

$ws.Range("A107").Value = "Estado de Palestina"
$ws.Range("B107").Value = 155
$ws.Range("C107").Value = 21
$ws.Range("D107").Value = 18
$ws.Range("E107").Value = 136
$ws.Range("F107").Value = 0
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 1

$ws.Range("A108").Value = "Sri Lanka"
$ws.Range("B108").Value = 148
$ws.Range("C108").Value = 2
$ws.Range("D108").Value = 21
$ws.Range("E108").Value = 124
$ws.Range("F108").Value = 5
$ws.Range("G108").Value = 0
$ws.Range("H108").Value = 3

$ws.Range("A109").Value = "Venezuela"
$ws.Range("B109").Value = 144
$ws.Range("C109").Value = 0
$ws.Range("D109").Value = 43
$ws.Range("E109").Value = 98
$ws.Range("F109").Value = 6
$ws.Range("G109").Value = 0
$ws.Range("H109").Value = 3

$ws.Range("A110").Value = "Montenegro"
$ws.Range("B110").Value = 140
$ws.Range("C110").Value = 17
$ws.Range("D110").Value = 0
$ws.Range("E110").Value = 138
$ws.Range("F110").Value = 4
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 2

$ws.Range("A111").Value = "Martinica"
$ws.Range("B111").Value = 135
$ws.Range("C111").Value = 0
$ws.Range("D111").Value = 27
$ws.Range("E111").Value = 105
$ws.Range("F111").Value = 16
$ws.Range("G111").Value = 0
$ws.Range("H111").Value = 3

$ws.Range("A125").Value = "Isla de Man"
$ws.Range("B125").Value = 75
$ws.Range("C125").Value = 7
$ws.Range("D125").Value = 0
$ws.Range("E125").Value = 74
$ws.Range("F125").Value = 0
$ws.Range("G125").Value = 0
$ws.Range("H125").Value = 1

$ws.Range("A126").Value = "Niger"
$ws.Range("B126").Value = 74
$ws.Range("C126").Value = 0
$ws.Range("D126").Value = 0
$ws.Range("E126").Value = 69
$ws.Range("F126").Value = 0
$ws.Range("G126").Value = 0
$ws.Range("H126").Value = 5

$ws.Range("A127").Value = "Liechtenstein"
$ws.Range("B127").Value = 72
$ws.Range("C127").Value = 0
$ws.Range("D127").Value = 0
$ws.Range("E127").Value = 72
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 0
